# Updated informations about languages
# My resume improved with my language skills

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("skills")

# The "Python" skill entry (column B, row 3) is merged with the more
# descriptive label already used in column G, so the plain "Python" shared
# string becomes unused and is dropped on save.
$ws.Range("B3").Value = "Python (system, parallel computing, programming)"

# Add language skills as new rows (French / English / Spanish), mirrored in
# both the A:B and F:G blocks like the rest of the "skills" sheet.
$ws.Range("A17").Value = "French"
$ws.Range("B17").Value = "C1/C2 (native tongue)"
$ws.Range("F17").Value = "French"
$ws.Range("G17").Value = "C1/C2 (native tongue)"

$ws.Range("A18").Value = "English "
$ws.Range("B18").Value = "B1/B2"
$ws.Range("F18").Value = "English "
$ws.Range("G18").Value = "B1/B2"

$ws.Range("A19").Value = "Spanish"
$ws.Range("B19").Value = "A2/B1"
$ws.Range("F19").Value = "Spanish"
$ws.Range("G19").Value = "A2/B1"

# The last language row (Spanish / A2/B1) carries a distinct font colour and
# wraps text, like pasted-in text would.
$ws.Range("B19").Font.Color = 4473924
$ws.Range("B19").WrapText = $true
$ws.Range("G19").Font.Color = 4473924
$ws.Range("G19").WrapText = $true

# Move the active tab/selection from "profil" to "skills", landing on F17.
$ws.Activate()
$ws.Range("F17").Select()
